$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1595.8235
$ws.Range("I28").Value = 1342
$ws.Range("J28").Value = 3499.5
$ws.Range("K28").Value = 1342
$ws.Range("L28").Value = 3499.5
$ws.Range("M28").Value = -857
$ws.Range("N28").Value = -4469.5

$ws.Range("H82").Value = 500
$ws.Range("I82").Value = 500
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 1500
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -1094

$ws.Range("H85").Value = 500
$ws.Range("I85").Value = 500
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 1500
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -96

$ws.Range("H131").Value = 3810.3157
$ws.Range("I131").Value = 1600.5
$ws.Range("J131").Value = 9997.799999999999
$ws.Range("K131").Value = 4801.5
$ws.Range("L131").Value = 29993.4
$ws.Range("M131").Value = 238.5
$ws.Range("N131").Value = -40073.39999999999

$ws.Range("H135").Value = 1023.1875
$ws.Range("I135").Value = 958.06665
$ws.Range("J135").Value = 2000
$ws.Range("K135").Value = 8622.599850000001
$ws.Range("L135").Value = 18000
$ws.Range("M135").Value = -6087.599850000001
$ws.Range("N135").Value = -23070

$ws.Range("H138").Value = 6443.457
$ws.Range("I138").Value = 6565.727
$ws.Range("J138").Value = 6387.4165
$ws.Range("K138").Value = 19697.181
$ws.Range("L138").Value = 19162.2495
$ws.Range("M138").Value = -14557.181
$ws.Range("N138").Value = -29442.2495

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1950
$ws.Range("I45").Value = 1950
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 1950
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -1573

$ws.Range("H61").Value = 6599.067
$ws.Range("I61").Value = 6602.769
$ws.Range("J61").Value = 6575
$ws.Range("K61").Value = 6602.769
$ws.Range("L61").Value = 6575
$ws.Range("M61").Value = -6390.769
$ws.Range("N61").Value = -6999

$ws.Range("H125").Value = 84995
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 84995
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 84995
$ws.Range("N125").Value = -94835

$ws.Range("H132").Value = 6898.3125
$ws.Range("I132").Value = 6490.9287
$ws.Range("J132").Value = 9750
$ws.Range("K132").Value = 19472.7861
$ws.Range("L132").Value = 29250
$ws.Range("M132").Value = -16942.7861
$ws.Range("N132").Value = -34310

$ws.Range("H136").Value = 6599.067
$ws.Range("I136").Value = 6602.769
$ws.Range("J136").Value = 6575
$ws.Range("K136").Value = 19808.307
$ws.Range("L136").Value = 19725
$ws.Range("M136").Value = -17258.307
$ws.Range("N136").Value = -24825

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 18580970
$ws.Range("I86").Value = 72747.62
$ws.Range("J86").Value = 66702348
$ws.Range("K86").Value = 72747.62
$ws.Range("L86").Value = 66702348
$ws.Range("M86").Value = -71624.62
$ws.Range("N86").Value = -66704594

$ws.Range("H89").Value = 18580970
$ws.Range("I89").Value = 72747.62
$ws.Range("J89").Value = 66702348
$ws.Range("K89").Value = 363738.1
$ws.Range("L89").Value = 333511740
$ws.Range("M89").Value = -358122.1
$ws.Range("N89").Value = -333522972

$ws.Range("H133").Value = 123497
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 123497
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 123497
$ws.Range("N133").Value = -133617

$ws.Range("H134").Value = 1086.625
$ws.Range("I134").Value = 956.1429000000001
$ws.Range("J134").Value = 2000
$ws.Range("K134").Value = 2868.4287
$ws.Range("L134").Value = 6000
$ws.Range("M134").Value = -333.4287000000004
$ws.Range("N134").Value = -11070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3146.1177
$ws.Range("I31").Value = 2748.8572
$ws.Range("J31").Value = 5000
$ws.Range("K31").Value = 2748.8572
$ws.Range("L31").Value = 5000
$ws.Range("M31").Value = -2453.8572
$ws.Range("N31").Value = -5590

$ws.Range("H34").Value = 3146.1177
$ws.Range("I34").Value = 2748.8572
$ws.Range("J34").Value = 5000
$ws.Range("K34").Value = 2748.8572
$ws.Range("L34").Value = 5000
$ws.Range("M34").Value = -2546.8572
$ws.Range("N34").Value = -5404

$ws.Range("H122").Value = 4616.5454
$ws.Range("I122").Value = 5523.375
$ws.Range("J122").Value = 2198.3333
$ws.Range("K122").Value = 16570.125
$ws.Range("L122").Value = 6594.999899999999
$ws.Range("M122").Value = -14120.125
$ws.Range("N122").Value = -11494.9999

$ws.Range("H134").Value = 2373.087
$ws.Range("I134").Value = 2173.0667
$ws.Range("J134").Value = 2748.125
$ws.Range("K134").Value = 6519.2001
$ws.Range("L134").Value = 8244.375
$ws.Range("M134").Value = -3984.2001
$ws.Range("N134").Value = -13314.375

$ws.Range("H140").Value = 137495
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 137495
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 137495
$ws.Range("N140").Value = -147855

$ws.Range("H141").Value = 659514.1
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 659514.1
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 659514.1
$ws.Range("N141").Value = -669874.1

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H112").Value = 9532.625
$ws.Range("I112").Value = 4865
$ws.Range("J112").Value = 11654.272
$ws.Range("K112").Value = 14595
$ws.Range("L112").Value = 34962.81600000001
$ws.Range("M112").Value = -13487
$ws.Range("N112").Value = -37178.81600000001

$ws.Range("H139").Value = 4171.125
$ws.Range("I139").Value = 2648
$ws.Range("J139").Value = 14833
$ws.Range("K139").Value = 7944
$ws.Range("L139").Value = 44499
$ws.Range("M139").Value = -2804
$ws.Range("N139").Value = -54779

$ws.Range("H140").Value = 1673000
$ws.Range("I140").Value = 1673000
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 5019000
$ws.Range("L140").Value = 0
$ws.Range("M140").Value = -5013820

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H63").Value = 75000
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 75000
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 75000
$ws.Range("N63").Value = -76372

$ws.Range("H66").Value = 75000
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 75000
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 225000
$ws.Range("N66").Value = -231864

$ws.Range("H132").Value = 30307018
$ws.Range("I132").Value = 4561.4443
$ws.Range("J132").Value = 166668060
$ws.Range("K132").Value = 13684.3329
$ws.Range("L132").Value = 500004180
$ws.Range("M132").Value = -11154.3329
$ws.Range("N132").Value = -500009240

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1385.7
$ws.Range("I16").Value = 1384.1111
$ws.Range("J16").Value = 1400
$ws.Range("K16").Value = 1384.1111
$ws.Range("L16").Value = 1400
$ws.Range("M16").Value = -1214.1111
$ws.Range("N16").Value = -1740

$ws.Range("H22").Value = 33333998
$ws.Range("I22").Value = 644.6667
$ws.Range("J22").Value = 55556230
$ws.Range("K22").Value = 644.6667
$ws.Range("L22").Value = 55556230
$ws.Range("M22").Value = -349.6667
$ws.Range("N22").Value = -55556820

$ws.Range("H27").Value = 33333998
$ws.Range("I27").Value = 644.6667
$ws.Range("J27").Value = 55556230
$ws.Range("K27").Value = 644.6667
$ws.Range("L27").Value = 55556230
$ws.Range("M27").Value = -537.6667
$ws.Range("N27").Value = -55556444

$ws.Range("H46").Value = 1865.2
$ws.Range("I46").Value = 1099.75
$ws.Range("J46").Value = 2375.5
$ws.Range("K46").Value = 1099.75
$ws.Range("L46").Value = 2375.5
$ws.Range("M46").Value = -911.75
$ws.Range("N46").Value = -2751.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2129.3809
$ws.Range("I122").Value = 2206.3684
$ws.Range("J122").Value = 1398
$ws.Range("K122").Value = 6619.1052
$ws.Range("L122").Value = 4194
$ws.Range("M122").Value = -4169.1052
$ws.Range("N122").Value = -9094

$ws.Range("H132").Value = 55556424
$ws.Range("I132").Value = 919.4706
$ws.Range("J132").Value = 1000000000
$ws.Range("K132").Value = 2758.4118
$ws.Range("L132").Value = 3000000000
$ws.Range("M132").Value = -228.4117999999999
$ws.Range("N132").Value = -3000005060
